$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.404.57"
$ws.Range("E2").Value = "  +3.85%  "
$ws.Range("D3").Value = "3.252.68"
$ws.Range("E3").Value = "  +3.65%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'577.92"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "'181.59"
$ws.Range("E6").Value = "  +7.87%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("E8").Value = "  -2.30%  "
$ws.Range("D9").Value = "3.251.80"
$ws.Range("E9").Value = "  +3.70%  "
$ws.Range("D10").Value = "'0.131"
$ws.Range("E10").Value = "  +6.39%  "
$ws.Range("E11").Value = "  +3.08%  "
$ws.Range("D12").Value = "'0.415"
$ws.Range("E12").Value = "  +6.52%  "
$ws.Range("D13").Value = "3.813.98"
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "'28.42"
$ws.Range("E15").Value = "  +5.79%  "
$ws.Range("D16").Value = "67.405.37"
$ws.Range("E16").Value = "  +3.99%  "
$ws.Range("D17").Value = "'0.0000167"
$ws.Range("E17").Value = "  +3.96%  "
$ws.Range("D18").Value = "3.255.74"
$ws.Range("E18").Value = "  +3.29%  "
$ws.Range("D19").Value = "'5.85"
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("D20").Value = "'13.57"
$ws.Range("E20").Value = "  +6.57%  "
$ws.Range("D21").Value = "'375.38"
$ws.Range("E21").Value = "  +5.54%  "
$ws.Range("E22").Value = "  +5.67%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'71.12"
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("E25").Value = "  +3.71%  "
$ws.Range("E26").Value = "  +3.85%  "
$ws.Range("D27").Value = "'9.56"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").Value = "'5.73"
$ws.Range("E30").Value = "  +8.79%  "
$ws.Range("E31").Value = "  +3.66%  "
$ws.Range("D32").Value = "'22.69"
$ws.Range("E32").Value = "  +3.98%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  +6.38%  "
$ws.Range("D35").Value = "'6.92"
$ws.Range("E35").Value = "  +5.27%  "
$ws.Range("D36").Value = "'163.53"
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("E37").Value = "  +5.04%  "
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("D39").Value = "'1.85"
$ws.Range("E39").Value = "  +4.97%  "
$ws.Range("E40").Value = "  +13.04%  "
$ws.Range("D41").Value = "'4.66"
$ws.Range("E41").Value = "  +12.45%  "
$ws.Range("D42").Value = "'26.64"
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("E43").Value = "  +7.69%  "
$ws.Range("D44").Value = "'356.44"
$ws.Range("E44").Value = "  +11.69%  "
$ws.Range("D45").Value = "2.704.21"
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("D46").Value = "'25.42"
$ws.Range("E46").Value = "  +6.49%  "
$ws.Range("D47").Value = "'40.83"
$ws.Range("E47").Value = "  +3.58%  "
$ws.Range("E48").Value = "  +4.68%  "
$ws.Range("E49").Value = "  +3.47%  "
$ws.Range("E50").Value = "  +7.07%  "
$ws.Range("E51").Value = "  -0.50%  "
